# Annex1 "more forms are added" edit
#
# 1. The title FORMTEXT field ("Text8") result is corrected:
#      "Eerste-Principes Studies van Nieuwe Tweedimensionale Materialen en Hun Fysische Eigenschappen"
#    becomes
#      "First-principes Studies van Nieuwe Tweedimensionale Materialen en hun Fysische Eigenschappen"
#    i.e. "Eerste-Principes" -> "First-principes" and "Hun" -> "hun".
#
# 2. The user's cursor (the "_GoBack" bookmark Word drops at the last edit
#    point) ends up right after "...Fysische Eig", inside that same field
#    result.
#
# 3. The date FORMTEXT field ("Text19") result was previously split across
#    two runs ("27/" and "11/2017") with the old "_GoBack" bookmark sitting
#    between them. That bookmark is gone now (it moved to the title, see
#    above) and the two runs collapse back into a single "27/11/2017" run.

$d = $word.ActiveDocument

# --- 1. Fix the title text -------------------------------------------------
$d.Content.Find.Execute("Eerste-Principes", $true, $false, $false, $false, $false, `
    $true, 1, $false, "First-principes", 2)

$d.Content.Find.Execute("Hun Fysische", $true, $false, $false, $false, $false, `
    $true, 1, $false, "hun Fysische", 2)

# --- 2. Collapse the date field's two runs back into one -------------------
# (This also removes the "_GoBack" bookmark that used to live between them,
#  since Word only ever keeps a single "_GoBack" bookmark in the document.)
$d.Content.Find.Execute("27/11/2017", $true, $false, $false, $false, $false, `
    $true, 1, $false, "27/11/2017", 2)

# --- 3. Drop the "_GoBack" bookmark at the new cursor location -------------
# (right after "...Fysische Eig", before "enschappen", inside the title
#  field's result text)
$titleRange = $d.Content
$titleRange.Find.Execute("Fysische Eig", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0)
$goBackPoint = $d.Range($titleRange.End, $titleRange.End)
$d.Bookmarks.Add("_GoBack", $goBackPoint)
